$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4617176651954651
$ws.Range("B1").Value = 0.5851160883903503
$ws.Range("C1").Value = 0.8609979748725891
$ws.Range("D1").Value = 3.997616529464722
$ws.Range("E1").Value = 4.176104545593262
